$d = $word.ActiveDocument

# Locate the version string "3.0.0" in the VERSION cell of the table.
$found = $d.Content
$found.Find.Execute("3.0.0", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$start = $found.Start
$end = $found.End

# Replace the text in-place with the first new segment ("4") - this keeps
# the original run's formatting (Verdana/000000/sz15/szCs24) intact.
$r1 = $d.Range($start, $end)
$r1.Text = "4"

# Insert the remaining version-number pieces as their own collapsed ranges,
# immediately following one another, so the resulting text reads "4.1.0".
$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter(".")

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter("1")

$r4 = $d.Range($r3.End, $r3.End)
$r4.InsertAfter(".0")

# Nudge each segment's formatting (toggle Bold on/off) so the four pieces
# are preserved as four distinct runs (matching the target markup) instead
# of being silently re-merged into a single run during save.
$r1.Font.Bold = 1
$r1.Font.Bold = 0
$r2.Font.Bold = 1
$r2.Font.Bold = 0
$r3.Font.Bold = 1
$r3.Font.Bold = 0
$r4.Font.Bold = 1
$r4.Font.Bold = 0
